$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text renames: humoment -> shape
$ws.Range("J1").Value = "shape"
$ws.Range("R1").Value = "texture-shape"
$ws.Range("V1").Value = "color-shape"
$ws.Range("Z1").Value = "texture-color-shape"

# Updated numeric results for rows 4-7
$ws.Range("B4").Value = 0.5610378380419134
$ws.Range("C4").Value = 0.583
$ws.Range("D4").Value = 0.5647766200071546
$ws.Range("E4").Value = 0.5839999999999999
$ws.Range("F4").Value = 0.596471664039624
$ws.Range("G4").Value = 0.628
$ws.Range("H4").Value = 0.6129084271147054
$ws.Range("I4").Value = 0.6020000000000001
$ws.Range("J4").Value = 0.6738780318888157
$ws.Range("K4").Value = 0.969
$ws.Range("L4").Value = 0.517122287669734
$ws.Range("M4").Value = 0.5309999999999999
$ws.Range("N4").Value = 0.6254024427147131
$ws.Range("O4").Value = 0.645
$ws.Range("P4").Value = 0.6437845128202339
$ws.Range("Q4").Value = 0.63
$ws.Range("R4").Value = 0.5646288753514351
$ws.Range("S4").Value = 0.588
$ws.Range("T4").Value = 0.5665949757264891
$ws.Range("U4").Value = 0.5855
$ws.Range("V4").Value = 0.5912022352494647
$ws.Range("W4").Value = 0.62
$ws.Range("X4").Value = 0.6114362808504905
$ws.Range("Y4").Value = 0.5995
$ws.Range("Z4").Value = 0.6266629720922198
$ws.Range("AA4").Value = 0.646
$ws.Range("AB4").Value = 0.6443363136392586
$ws.Range("AC4").Value = 0.631
$ws.Range("B5").Value = 0.5810324902546983
$ws.Range("C5").Value = 0.605
$ws.Range("D5").Value = 0.5810774716889557
$ws.Range("E5").Value = 0.5905
$ws.Range("F5").Value = 0.7116478961297652
$ws.Range("G5").Value = 0.779
$ws.Range("H5").Value = 0.6721240867499313
$ws.Range("I5").Value = 0.6839999999999999
$ws.Range("J5").Value = 0.6206504386321211
$ws.Range("K5").Value = 0.764
$ws.Range("L5").Value = 0.5379334277740289
$ws.Range("M5").Value = 0.5545
$ws.Range("N5").Value = 0.6384016540939585
$ws.Range("O5").Value = 0.6570000000000001
$ws.Range("P5").Value = 0.6485119546252436
$ws.Range("Q5").Value = 0.636
$ws.Range("R5").Value = 0.5830725719375816
$ws.Range("S5").Value = 0.608
$ws.Range("T5").Value = 0.5829378049724532
$ws.Range("U5").Value = 0.592
$ws.Range("V5").Value = 0.7016975702337404
$ws.Range("W5").Value = 0.773
$ws.Range("X5").Value = 0.6590904839159274
$ws.Range("Y5").Value = 0.6725000000000001
$ws.Range("Z5").Value = 0.6363280508781239
$ws.Range("AA5").Value = 0.6570000000000001
$ws.Range("AB5").Value = 0.6457375535616043
$ws.Range("AC5").Value = 0.633
$ws.Range("B6").Value = 0.5542377678943767
$ws.Range("C6").Value = 0.5569999999999999
$ws.Range("D6").Value = 0.5824990502093139
$ws.Range("E6").Value = 0.5974999999999999
$ws.Range("F6").Value = 0.6890143946505687
$ws.Range("G6").Value = 0.7610000000000001
$ws.Range("H6").Value = 0.6478327887034013
$ws.Range("I6").Value = 0.6624999999999999
$ws.Range("J6").Value = 0.6456260945889024
$ws.Range("K6").Value = 0.857
$ws.Range("L6").Value = 0.5249892209954442
$ws.Range("M6").Value = 0.5389999999999999
$ws.Range("N6").Value = 0.6428852974078167
$ws.Range("O6").Value = 0.659
$ws.Range("P6").Value = 0.6547624066178906
$ws.Range("Q6").Value = 0.6445000000000001
$ws.Range("R6").Value = 0.5576126680829809
$ws.Range("S6").Value = 0.5589999999999999
$ws.Range("T6").Value = 0.5863015526479983
$ws.Range("U6").Value = 0.599
$ws.Range("V6").Value = 0.6694819440469891
$ws.Range("W6").Value = 0.748
$ws.Range("X6").Value = 0.6228876461889888
$ws.Range("Y6").Value = 0.6345000000000001
$ws.Range("Z6").Value = 0.6348496035251657
$ws.Range("AA6").Value = 0.6489999999999999
$ws.Range("AB6").Value = 0.6502375157952499
$ws.Range("AC6").Value = 0.639
$ws.Range("B7").Value = 0.5064608506989259
$ws.Range("C7").Value = 0.5409999999999999
$ws.Range("D7").Value = 0.4837156620475147
$ws.Range("E7").Value = 0.496
$ws.Range("F7").Value = 0.4748255085823812
$ws.Range("G7").Value = 0.474
$ws.Range("H7").Value = 0.4931001515121298
$ws.Range("I7").Value = 0.4970000000000001
$ws.Range("J7").Value = 0.6018408147846502
$ws.Range("K7").Value = 0.8029999999999999
$ws.Range("L7").Value = 0.4894977347593315
$ws.Range("M7").Value = 0.5175
$ws.Range("N7").Value = 0.5679917654016056
$ws.Range("O7").Value = 0.592
$ws.Range("P7").Value = 0.551254268987578
$ws.Range("Q7").Value = 0.5605
$ws.Range("R7").Value = 0.5072526266468894
$ws.Range("S7").Value = 0.5359999999999999
$ws.Range("T7").Value = 0.4870293127680508
$ws.Range("U7").Value = 0.4955
$ws.Range("V7").Value = 0.5128940962268628
$ws.Range("W7").Value = 0.529
$ws.Range("X7").Value = 0.5072530046965837
$ws.Range("Y7").Value = 0.519
$ws.Range("Z7").Value = 0.5655732188573166
$ws.Range("AA7").Value = 0.583
$ws.Range("AB7").Value = 0.5544619749694173
$ws.Range("AC7").Value = 0.5605
